$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the contents of C2:O2 (values removed, keeping style/formatting)
$ws.Range("C2:O2").ClearContents()

# Update the active cell selection to K9
$ws.Range("K9").Select()
